$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$tr = $s.Shapes.Item(1).TextFrame.TextRange
$tr.Text = "Slide 3 title was "
$tr.InsertAfter("also changed")
